$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# The subtitle placeholder starts out empty (just an endParaRPr). Build the
# paragraph up run-by-run with InsertAfter so the trailing endParaRPr (and
# its formatting) is preserved, exactly like typing the text in the UI.
$run1 = $tr.InsertAfter("This is ")
$run1.LanguageID = "en-US"

$run2 = $tr.InsertAfter("Fablehaft")
$run2.LanguageID = "en-US"
